$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Step_1")
$ws.Range("A71").Value = "dw00"
$ws.Range("B71").Value = "dw10"
